$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "1"
